# Updates the cryptos list (prices / 1h volume changes) to the latest
# scraped values, and fixes the ranking order of a few coins whose
# positions shifted (rows 45-48): FraxShare, HuobiToken, ARBITRUM, Cronos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.902.88"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.128.25"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.18"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.34"
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.94"
$ws.Range("E12").Value = "  +6.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.442.03"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.13"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.810"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.142.71"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.900.52"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.82"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.73"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.67"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.47"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.48"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  +9.49%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("E34").Value = "  +11.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.80"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.17"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.74"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.529.91"
$ws.Range("E44").Value = "  +6.52%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.84"
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0917"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.326.98"
$ws.Range("E51").Value = "  +1.66%  "
